# Natmi following Dr Hou advice
# Update ligand/receptor expressing-cell counts (E, K: 1 -> 3) and
# the dependent NATMI-derived statistics for Tgfb3-Tgfbr1 LR pairs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.277420333333333
$ws.Range("H2").Value = 3.832261
$ws.Range("I2").Value = 0.01913942624337554
$ws.Range("J2").Value = 0.01913942624337554
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 44.50020533333333
$ws.Range("N2").Value = 133.500616
$ws.Range("O2").Value = 0.2926972930209797
$ws.Range("P2").Value = 0.2926972930209797
$ws.Range("Q2").Value = 56.84546713030845
$ws.Range("R2").Value = 511.609204172776
$ws.Range("S2").Value = 0.00560205825141072
$ws.Range("T2").Value = 0.00560205825141072

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.277420333333333
$ws.Range("H3").Value = 3.832261
$ws.Range("I3").Value = 0.01913942624337554
$ws.Range("J3").Value = 0.01913942624337554
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 28.185334
$ws.Range("N3").Value = 84.55600199999999
$ws.Range("O3").Value = 0.1853872561462678
$ws.Range("P3").Value = 0.1853872561462678
$ws.Range("Q3").Value = 36.00451875339133
$ws.Range("R3").Value = 324.040668780522
$ws.Range("S3").Value = 0.003548205715473262
$ws.Range("T3").Value = 0.003548205715473262

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.277420333333333
$ws.Range("H4").Value = 3.832261
$ws.Range("I4").Value = 0.01913942624337554
$ws.Range("J4").Value = 0.01913942624337554
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 62.31760433333333
$ws.Range("N4").Value = 186.952813
$ws.Range("O4").Value = 0.4098901108273345
$ws.Range("P4").Value = 0.4098901108273344
$ws.Range("Q4").Value = 79.60577490002144
$ws.Range("R4").Value = 716.451974100193
$ws.Range("S4").Value = 0.007845061544068796
$ws.Range("T4").Value = 0.007845061544068794

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.277420333333333
$ws.Range("H5").Value = 3.832261
$ws.Range("I5").Value = 0.01913942624337554
$ws.Range("J5").Value = 0.01913942624337554
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.031762
$ws.Range("N5").Value = 51.09528599999999
$ws.Range("O5").Value = 0.1120253400054181
$ws.Range("P5").Value = 0.1120253400054181
$ws.Range("Q5").Value = 21.756719091294
$ws.Range("R5").Value = 195.810471821646
$ws.Range("S5").Value = 0.002144100732422767
$ws.Range("T5").Value = 0.002144100732422767

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 45.44725166666667
$ws.Range("H6").Value = 136.341755
$ws.Range("I6").Value = 0.6809303864519871
$ws.Range("J6").Value = 0.6809303864519872
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 44.50020533333333
$ws.Range("N6").Value = 133.500616
$ws.Range("O6").Value = 0.2926972930209797
$ws.Range("P6").Value = 0.2926972930209797
$ws.Range("Q6").Value = 2022.412031002342
$ws.Range("R6").Value = 18201.70827902108
$ws.Range("S6").Value = 0.1993064808502262
$ws.Range("T6").Value = 0.1993064808502262

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 45.44725166666667
$ws.Range("H7").Value = 136.341755
$ws.Range("I7").Value = 0.6809303864519871
$ws.Range("J7").Value = 0.6809303864519872
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 28.185334
$ws.Range("N7").Value = 84.55600199999999
$ws.Range("O7").Value = 0.1853872561462678
$ws.Range("P7").Value = 0.1853872561462678
$ws.Range("Q7").Value = 1280.945967607056
$ws.Range("R7").Value = 11528.51370846351
$ws.Range("S7").Value = 0.1262358159709517
$ws.Range("T7").Value = 0.1262358159709517

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 45.44725166666667
$ws.Range("H8").Value = 136.341755
$ws.Range("I8").Value = 0.6809303864519871
$ws.Range("J8").Value = 0.6809303864519872
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 62.31760433333333
$ws.Range("N8").Value = 186.952813
$ws.Range("O8").Value = 0.4098901108273345
$ws.Range("P8").Value = 0.4098901108273344
$ws.Range("Q8").Value = 2832.163847400757
$ws.Range("R8").Value = 25489.47462660681
$ws.Range("S8").Value = 0.2791066315685047
$ws.Range("T8").Value = 0.2791066315685047

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 45.44725166666667
$ws.Range("H9").Value = 136.341755
$ws.Range("I9").Value = 0.6809303864519871
$ws.Range("J9").Value = 0.6809303864519872
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 17.031762
$ws.Range("N9").Value = 51.09528599999999
$ws.Range("O9").Value = 0.1120253400054181
$ws.Range("P9").Value = 0.1120253400054181
$ws.Range("Q9").Value = 774.0467739407699
$ws.Range("R9").Value = 6966.420965466929
$ws.Range("S9").Value = 0.0762814580623046
$ws.Range("T9").Value = 0.0762814580623046

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.4966396666666666
$ws.Range("H10").Value = 1.489919
$ws.Range("I10").Value = 0.007441088905245192
$ws.Range("J10").Value = 0.007441088905245193
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 44.50020533333333
$ws.Range("N10").Value = 133.500616
$ws.Range("O10").Value = 0.2926972930209797
$ws.Range("P10").Value = 0.2926972930209797
$ws.Range("Q10").Value = 22.10056714334489
$ws.Range("R10").Value = 198.905104290104
$ws.Range("S10").Value = 0.002177986579693713
$ws.Range("T10").Value = 0.002177986579693713

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.4966396666666666
$ws.Range("H11").Value = 1.489919
$ws.Range("I11").Value = 0.007441088905245192
$ws.Range("J11").Value = 0.007441088905245193
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 28.185334
$ws.Range("N11").Value = 84.55600199999999
$ws.Range("O11").Value = 0.1853872561462678
$ws.Range("P11").Value = 0.1853872561462678
$ws.Range("Q11").Value = 13.99795488264867
$ws.Range("R11").Value = 125.981593943838
$ws.Range("S11").Value = 0.001379483054883842
$ws.Range("T11").Value = 0.001379483054883842

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.4966396666666666
$ws.Range("H12").Value = 1.489919
$ws.Range("I12").Value = 0.007441088905245192
$ws.Range("J12").Value = 0.007441088905245193
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 62.31760433333333
$ws.Range("N12").Value = 186.952813
$ws.Range("O12").Value = 0.4098901108273345
$ws.Range("P12").Value = 0.4098901108273344
$ws.Range("Q12").Value = 30.94939424357188
$ws.Range("R12").Value = 278.544548192147
$ws.Range("S12").Value = 0.003050028756047
$ws.Range("T12").Value = 0.003050028756047001

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.4966396666666666
$ws.Range("H13").Value = 1.489919
$ws.Range("I13").Value = 0.007441088905245192
$ws.Range("J13").Value = 0.007441088905245193
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 17.031762
$ws.Range("N13").Value = 51.09528599999999
$ws.Range("O13").Value = 0.1120253400054181
$ws.Range("P13").Value = 0.1120253400054181
$ws.Range("Q13").Value = 8.458648602425997
$ws.Range("R13").Value = 76.12783742183399
$ws.Range("S13").Value = 0.000833590514620637
$ws.Range("T13").Value = 0.000833590514620637

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 19.52156333333333
$ws.Range("H14").Value = 58.56469
$ws.Range("I14").Value = 0.2924890983993922
$ws.Range("J14").Value = 0.2924890983993922
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 44.50020533333333
$ws.Range("N14").Value = 133.500616
$ws.Range("O14").Value = 0.2926972930209797
$ws.Range("P14").Value = 0.2926972930209797
$ws.Range("Q14").Value = 868.7135767610044
$ws.Range("R14").Value = 7818.42219084904
$ws.Range("S14").Value = 0.08561076733964906
$ws.Range("T14").Value = 0.08561076733964908

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 19.52156333333333
$ws.Range("H15").Value = 58.56469
$ws.Range("I15").Value = 0.2924890983993922
$ws.Range("J15").Value = 0.2924890983993922
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 28.185334
$ws.Range("N15").Value = 84.55600199999999
$ws.Range("O15").Value = 0.1853872561462678
$ws.Range("P15").Value = 0.1853872561462678
$ws.Range("Q15").Value = 550.2217827521532
$ws.Range("R15").Value = 4951.996044769379
$ws.Range("S15").Value = 0.05422375140495905
$ws.Range("T15").Value = 0.05422375140495906

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 19.52156333333333
$ws.Range("H16").Value = 58.56469
$ws.Range("I16").Value = 0.2924890983993922
$ws.Range("J16").Value = 0.2924890983993922
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 62.31760433333333
$ws.Range("N16").Value = 186.952813
$ws.Range("O16").Value = 0.4098901108273345
$ws.Range("P16").Value = 0.4098901108273344
$ws.Range("Q16").Value = 1216.537059774774
$ws.Range("R16").Value = 10948.83353797297
$ws.Range("S16").Value = 0.119888388958714
$ws.Range("T16").Value = 0.119888388958714

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 19.52156333333333
$ws.Range("H17").Value = 58.56469
$ws.Range("I17").Value = 0.2924890983993922
$ws.Range("J17").Value = 0.2924890983993922
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 17.031762
$ws.Range("N17").Value = 51.09528599999999
$ws.Range("O17").Value = 0.1120253400054181
$ws.Range("P17").Value = 0.1120253400054181
$ws.Range("Q17").Value = 332.4866205612599
$ws.Range("R17").Value = 2992.37958505134
$ws.Range("S17").Value = 0.0327661906960701
$ws.Range("T17").Value = 0.0327661906960701
